# "nota EV02 y sprint 2"
#
# - B7: drop "en PDF " from the boletas/ventas user story text.
# - B8: reword the "boleta digital" story to "comprobante digital" /
#        "constancia" instead of "comprobante".
# - Move the active selection from C9 to A20 (and scroll so row 8 is
#   the first visible row, matching the saved view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "Como cliente quiero recibir un comprobante digital para tener constancia de mi compra."
$ws.Range("B7").Value = "Como dueña de la florería necesito registrar ventas y generar boletas para entregar comprobantes a clientes."

$ws.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("A20").Select()
